# Applies the OOXML diff:
#   1) Three runs currently highlighted "red" become highlighted "green"
#      (wdBrightGreen = 4, which serialises as w:highlight w:val="green").
#   2) A block of runs currently colored FF0000 (red) becomes 00B050
#      (a shade of green), including the paragraph-mark run properties
#      of the "Sistema de perfil de usuarios" paragraph.

$d = $word.ActiveDocument

# --- 1a) " de los 10 mejor ranqueados o aquellos que fueron indicados como"
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = " de los 10 mejor ranqueados o aquellos que fueron indicados como"
$found = $find.Execute()
if ($found) {
    $rng.Font.HighlightColorIndex = 4
}

# --- 1b) the space + "Ranqueados: Mayor a menor " pair of runs
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = " Ranqueados: Mayor a menor "
$found = $find.Execute()
if ($found) {
    $rng.Font.HighlightColorIndex = 4
}

# --- 2a) "● Sistema de perfil de usuarios " paragraph: run color AND
#         paragraph-mark run color (FF0000 -> 00B050)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*● Sistema de perfil de usuarios*") {
        $p.Range.Font.Color = 5287936
    }
}

# --- 2b) "Al dar de alta un usuario se deben indicar el perfil " /
#         "del mismo" / ". " runs (FF0000 -> 00B050)
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "Al dar de alta un usuario se deben indicar el perfil del mismo. "
$found = $find.Execute()
if ($found) {
    $rng.Font.Color = 5287936
}
